# Adapt tests to control version
#
# The "settings" sheet of this ODK-style form workbook currently has two
# columns (form_title, form_id). We add a third "version" column so the
# test fixture can control/assert the form version, giving it the value 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Make sure we're editing/selecting on the right sheet (it's already the
# active tab in this workbook, but be explicit).
$ws.Activate()

# New header cell for the version column, next to the existing
# form_title/form_id headers.
$ws.Range("C1").Value = "version"

# New data cell: the actual version number for this form.
$ws.Range("C2").Value = 1

# Mirror the cursor ending up one row below the newly entered value, as if
# the user had just typed it in and pressed Enter.
$ws.Range("C3").Select()
